$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.85"
$ws.Range("E2").Value = "'0.99%"
$ws.Range("D3").Value = "'32.17"
$ws.Range("E3").Value = "'1.11%"
$ws.Range("D4").Value = "'4.995"
$ws.Range("E4").Value = "'-2.90%"
$ws.Range("D5").Value = "'0.07912"
$ws.Range("E5").Value = "'-3.38%"
$ws.Range("D6").Value = "'2.110"
$ws.Range("E6").Value = "'-17.28%"
$ws.Range("D7").Value = "'7.865"
$ws.Range("E7").Value = "'0.17%"
$ws.Range("D8").Value = "'3.809"
$ws.Range("E8").Value = "'-1.21%"
$ws.Range("D9").Value = "'0.9266"
$ws.Range("E9").Value = "'-0.04%"
$ws.Range("D10").Value = "'0.1756"
$ws.Range("E10").Value = "'-0.13%"
$ws.Range("D11").Value = "'0.08039"
$ws.Range("E11").Value = "'7.38%"
$ws.Range("D12").Value = "'0.08755"
$ws.Range("E12").Value = "'-2.36%"
$ws.Range("D13").Value = "'0.03159"
$ws.Range("E13").Value = "'4.38%"
$ws.Range("D14").Value = "'0.1004"
$ws.Range("E14").Value = "'0.10%"
$ws.Range("D15").Value = "'0.001513"
$ws.Range("E15").Value = "'0.14%"
$ws.Range("D16").Value = "'0.006004"
$ws.Range("E16").Value = "'-0.62%"
$ws.Range("D17").Value = "'3.468"
$ws.Range("E17").Value = "'-3.76%"
$ws.Range("D18").Value = "'2.279"
$ws.Range("E18").Value = "'-0.26%"
$ws.Range("D20").Value = "'0.1290"
$ws.Range("E20").Value = "'-4.20%"
$ws.Range("D21").Value = "'4.189"
$ws.Range("E21").Value = "'-1.33%"
$ws.Range("E22").Value = "'6.75%"
$ws.Range("D23").Value = "'0.04607"
$ws.Range("E23").Value = "'-0.59%"
$ws.Range("D24").Value = "'0.001237"
$ws.Range("E24").Value = "'-0.91%"
$ws.Range("D25").Value = "'0.004493"
$ws.Range("E25").Value = "'-1.27%"
$ws.Range("E26").Value = "'4.09%"
$ws.Range("D39").Value = "'0.01738"
$ws.Range("E39").Value = "'-2.53%"
$ws.Range("D40").Value = "'0.04805"
$ws.Range("E40").Value = "'4.40%"
$ws.Range("D41").Value = "'0.007330"
$ws.Range("E41").Value = "'6.78%"
$ws.Range("D42").Value = "'0.1367"
$ws.Range("E42").Value = "'-1.02%"
$ws.Range("D43").Value = "'0.002360"
$ws.Range("E43").Value = "'10.17%"
$ws.Range("D44").Value = "'0.01108"
$ws.Range("E44").Value = "'12.32%"
$ws.Range("D45").Value = "'0.00006048"
$ws.Range("E45").Value = "'-2.35%"
$ws.Range("E46").Value = "'-0.09%"
$ws.Range("E47").Value = "'-59.60%"
$ws.Range("D48").Value = "'0.8234"
$ws.Range("E48").Value = "'2.26%"
$ws.Range("E49").Value = "'-0.09%"
$ws.Range("E50").Value = "'-0.09%"
